# Actualización automática hashcode jue mar 12 01:23:38 CET 2020
# Updates the "hashcode" value column (B) for a set of rows identified by
# their existing (old) hash values, replacing them with newly computed ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B52"  = "3eec821c25c09918accaedc26ad901d8"
    "B67"  = "db00c40cc339b31d2c635f5ad69fa42b"
    "B122" = "b99c222ca5e4a403aa2cd2aa985c2a97"
    "B142" = "2458a3c112b75bf87ad4ba14997740c8"
    "B145" = "967fd68316881401e27837c494381944"
    "B167" = "4d0d47f36c6cc070a68eff2592d6980f"
    "B230" = "0defa52e30cb2325fa7a2ac126c3721d"
    "B233" = "8248ebf504f034f12a563a9ea82aafed"
    "B331" = "6a21293b32a007b915df89a1aa014239"
    "B342" = "9e52312f5152c5e9ce4e4508b3e10763"
    "B343" = "968970307ea09d1b31218e3d61f00702"
    "B472" = "26f745cf167fbdd38fdd9c8f60253270"
    "B548" = "83c63fb7650d2edcaf66310fc00fd6ec"
    "B599" = "e551f3a2b05431cda4a45457e6e6d340"
    "B601" = "29ed2f7634f3ac085573cafce281d045"
    "B632" = "cc961aa8e1a2abd2c4952c818ab49123"
    "B681" = "656a006cbcf7ec2b05e7d214903ff9f2"
    "B686" = "e4b51c1dd3ebb322a18638e0957ae921"
    "B715" = "e9d235c08e4c1c94c33f933f84de5cef"
    "B787" = "a22661b820cb2e7ed848468e08c78eed"
    "B839" = "1b6ed3b3a92747c4d0e386b2f9f87036"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
